$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: exp_long  (category totals by year; updated raw dollar figures)
# ---------------------------------------------------------------------------
$expLong = $wb.Worksheets.Item("exp_long")
$expLong.Range("C721").Value = 1676.32     # 2021 Community Development
$expLong.Range("C743").Value = 92807.11    # 2021 Total
$expLong.Range("C752").Value = 1422.14     # 2022 Community Development
$expLong.Range("C770").Value = 1879.3      # 2022 Revenue
$expLong.Range("C774").Value = 101829.1    # 2022 Total

# ---------------------------------------------------------------------------
# Sheet: Table 1  (FY2022 expenditure categories, sorted by $ billions desc)
# Rows 10-13 re-sort because "Revenue" dropped from 2.3bn to 1.9bn, pushing
# it below "University Education", "Tollway" and "Debt Service".
# ---------------------------------------------------------------------------
$table1 = $wb.Worksheets.Item("Table 1")

$table1.Range("A10").Value = "University Education"
$table1.Range("B10").Value = 2.3
$table1.Range("C10").Value = 4.72
$table1.Range("D10").Value = 0.44

$table1.Range("A11").Value = "Tollway"
$table1.Range("B11").Value = 2.1
$table1.Range("C11").Value = 7.21
$table1.Range("D11").Value = 7.54

$table1.Range("A12").Value = "Debt Service"
$table1.Range("B12").Value = 2
$table1.Range("C12").Value = -0.83
$table1.Range("D12").Value = 6.11

$table1.Range("A13").Value = "Revenue"
$table1.Range("B13").Value = 1.9
$table1.Range("C13").Value = 11.96
$table1.Range("D13").Value = 6.43

# Row 17: Community Development % change / CAGR recompute
$table1.Range("C17").Value = -15.16
$table1.Range("D17").Value = 4.77

# Row 32: Total row
$table1.Range("B32").Value = 101.8
$table1.Range("C32").Value = 9.72
$table1.Range("D32").Value = 5.05

# ---------------------------------------------------------------------------
# Sheet: Table 4.b  (CAGR table: 1/2/3/5/10/24-year growth by category)
# ---------------------------------------------------------------------------
$table4b = $wb.Worksheets.Item("Table 4.b")

# Row 8: Community Development
$table4b.Range("B8").Value = -15.16
$table4b.Range("C8").Value = 51.43
$table4b.Range("D8").Value = 35.14
$table4b.Range("E8").Value = 16.98
$table4b.Range("F8").Value = 3.31
$table4b.Range("G8").Value = 4.77

# Row 26: Revenue
$table4b.Range("B26").Value = 11.96
$table4b.Range("C26").Value = 29.18
$table4b.Range("D26").Value = 46.38
$table4b.Range("E26").Value = 30.84
$table4b.Range("F26").Value = 14.11
$table4b.Range("G26").Value = 6.43

# Row 32: Total
$table4b.Range("B32").Value = 9.72
$table4b.Range("C32").Value = 11.73
$table4b.Range("D32").Value = 11.04
$table4b.Range("E32").Value = 7.27
$table4b.Range("F32").Value = 5.46
$table4b.Range("G32").Value = 5.05

# ---------------------------------------------------------------------------
# Sheet: year_totals  (yearly expenditures / revenues / fiscal gap)
# ---------------------------------------------------------------------------
$yearTotals = $wb.Worksheets.Item("year_totals")

$yearTotals.Range("B25").Value = 92807.10818869   # 2021 Expenditures
$yearTotals.Range("D25").Value = -1001            # 2021 Fiscal Gap

$yearTotals.Range("B26").Value = 101829.10187407  # 2022 Expenditures
$yearTotals.Range("C26").Value = 113021.56536341  # 2022 Revenues
$yearTotals.Range("D26").Value = 11192            # 2022 Fiscal Gap
